$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E35").Value = "Error: 📝 Summarizing and translating - Still failed after 3 attempts: Error code: 504 - {'code': 50501, 'message': 'Model service timeout. Please try again later.', 'data': None}"
$ws.Range("E38").Value = "Error: 🎙️ Transcribing with Whisper - Expecting value: line 1 column 1 (char 0)"
$ws.Range("E39").Value = "Error: 🗣️ Generating audio - Command '['ffmpeg', '-i', 'output/audio/tmp/337_0_temp.wav', '-filter:a', 'atempo=-0.007', '-y', 'output/audio/segs/337_0.wav']' returned non-zero exit status 222."
$ws.Range("E40").Value = "Error: 🗣️ Generating audio - Command '['ffmpeg', '-i', 'output/audio/tmp/337_0_temp.wav', '-filter:a', 'atempo=-0.007', '-y', 'output/audio/segs/337_0.wav']' returned non-zero exit status 222."
